$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Expected 24V PSU Load value" and "Expected 2nd 24V PSU Load value"
# test-data cells from 0.020 to 0.000 (kept as text, matching the existing
# text/quote-prefixed formatting of these cells).
$ws.Range("F8").Value = "'0.000"
$ws.Range("J8").Value = "'0.000"

# Move the active selection to I13, as captured in the saved workbook view.
$ws.Range("I13").Select()
